$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  D=0.8289408254599702; E=0.7252124645892352; F=0.5089463220675944; G=0.5981308411214953},
    @{Row=3;  D=0.8340792308967346; E=0.7263814616755794; F=0.5400927766732936; G=0.6195362979855569},
    @{Row=4;  D=0.8385546162771423; E=0.7407740774077408; F=0.5453943008614976; G=0.6282442748091602},
    @{Row=5;  D=0.8347422509530913; E=0.7310469314079422; F=0.536779324055666;  G=0.619029423003439},
    @{Row=6;  D=0.8320901707276646; E=0.7242314647377939; F=0.5308151093439364; G=0.6126195028680689},
    @{Row=7;  D=0.8339134758826454; E=0.7332106715731371; F=0.5281643472498343; G=0.6140215716486903},
    @{Row=8;  D=0.834410740924913;  E=0.7305605786618445; F=0.5354539430086149; G=0.6179732313575526},
    @{Row=9;  D=0.8324216807558429; E=0.7263636363636363; F=0.5294897282968853; G=0.6124952088922959},
    @{Row=10; D=0.8334162108403779; E=0.7262118491921006; F=0.5361166335321405; G=0.6168509340449866}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("E$r").Value = $entry.E
    $ws.Range("F$r").Value = $entry.F
    $ws.Range("G$r").Value = $entry.G
}

$wb.Save()
